$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.208.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.284.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.394"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +34.73%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.646"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.279.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.587"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000265"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.176"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.893.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.318.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.310.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "435.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000131"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.180"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +23.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "553.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.399"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "180.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.30%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.748"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.54%  "
